$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.483.71'
$ws.Range('E2').Value = '  -1.07%  '

$ws.Range('D3').Value = '1.912.01'
$ws.Range('E3').Value = '  -1.41%  '

$ws.Range('D4').Value = '''0.9994'
$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').Value = '''239.03'
$ws.Range('E5').Value = '  -1.54%  '

$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  -0.05%  '

$ws.Range('D7').Value = '''0.4780'
$ws.Range('E7').Value = '  -2.23%  '

$ws.Range('D8').Value = '''0.2848'
$ws.Range('E8').Value = '  -3.10%  '

$ws.Range('D9').Value = '''0.06709'
$ws.Range('E9').Value = '  -2.89%  '

$ws.Range('D10').Value = '''18.91'
$ws.Range('E10').Value = '  -2.54%  '

$ws.Range('D11').Value = '''102.30'
$ws.Range('E11').Value = '  -3.23%  '

$ws.Range('D12').Value = '''0.07704'
$ws.Range('E12').Value = '  -0.22%  '

$ws.Range('D13').Value = '1.918.44'
$ws.Range('E13').Value = '  -1.22%  '

$ws.Range('D14').Value = '''5.213'
$ws.Range('E14').Value = '  -2.62%  '

$ws.Range('D15').Value = '''0.6711'
$ws.Range('E15').Value = '  -3.95%  '

$ws.Range('D16').Value = '''271.34'
$ws.Range('E16').Value = '  -0.70%  '

$ws.Range('D17').Value = '30.513.64'
$ws.Range('E17').Value = '  -1.00%  '

$ws.Range('D18').Value = '''1.001'
$ws.Range('E18').Value = '  +0.01%  '

$ws.Range('D19').Value = '''0.000007478'
$ws.Range('E19').Value = '  -3.11%  '

$ws.Range('D20').Value = '''12.69'
$ws.Range('E20').Value = '  -3.07%  '

$ws.Range('D21').Value = '''5.427'
$ws.Range('E21').Value = '  -1.62%  '

$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = '''0.9995'
$ws.Range('E22').Value = '  -0.19%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '''6.306'
$ws.Range('E23').Value = '  -3.91%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '''9.391'
$ws.Range('E24').Value = '  -3.35%  '

$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''166.94'
$ws.Range('E25').Value = '  +0.12%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '''19.36'
$ws.Range('E26').Value = '  -1.23%  '

$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '''2.059'
$ws.Range('E27').Value = '  -4.75%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '''1.381'
$ws.Range('E28').Value = '  -0.68%  '

$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').Value = '''0.1004'
$ws.Range('E29').Value = '  -3.35%  '

$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '''4.629'
$ws.Range('E30').Value = '  +1.23%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.513'
$ws.Range('E31').Value = '  -2.61%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''4.231'
$ws.Range('E32').Value = '  -3.03%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.04735'
$ws.Range('E33').Value = '  -2.53%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '''0.7259'
$ws.Range('E34').Value = '  -3.94%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '''1.109'
$ws.Range('E35').Value = '  -3.98%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '''2.716'
$ws.Range('E36').Value = '  -0.67%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.01923'
$ws.Range('E37').Value = '  -3.57%  '

$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''2.611'
$ws.Range('E38').Value = '  -1.80%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '''6.267'
$ws.Range('E39').Value = '  -3.47%  '

$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '''74.87'
$ws.Range('E40').Value = '  -3.30%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '''1.968'
$ws.Range('E41').Value = '  -5.79%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''0.8583'
$ws.Range('E42').Value = '  -5.06%  '

$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '''105.13'
$ws.Range('E43').Value = '  -2.68%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '''0.4265'
$ws.Range('E44').Value = '  -3.23%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '''1.000'
$ws.Range('E45').Value = '  +0.14%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '''7.419'
$ws.Range('E46').Value = '  -4.52%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '''0.1201'
$ws.Range('E47').Value = '  -3.65%  '

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '''917.44'
$ws.Range('E48').Value = '  -7.35%  '

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '''34.80'
$ws.Range('E49').Value = '  -3.60%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''8.831'
$ws.Range('E50').Value = '  -5.03%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.05754'
$ws.Range('E51').Value = '  +0.38%  '
